$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Date" in D1
$ws.Range("D1").Value = "Date"

# Fix B3 and B4 values (swap: B3 3->2, B4 2->3)
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3

# Add new rows 5-21, following pattern A=row-1, B=row-1, C=3
for ($row = 5; $row -le 21; $row++) {
    $val = $row - 1
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = $val
    $ws.Cells.Item($row, 3).Value = 3
}

# Update selection to A2:D4 with active cell A2
$ws.Range("A2:D4").Select()
